# Apply the edit described by the diff:
#  - Add 5 new header columns (E:I) to every worksheet, labelled 165, 176, 219, 230, 255
#  - Populate the corresponding new data cells (E2:I4) on every worksheet
#  - Update the pre-existing B2:D4 values on Sheet2 and Sheet3 (Sheet1's B2:D4 values are unchanged)

$wb = $excel.ActiveWorkbook

$headerCols = @("E", "F", "G", "H", "I")
$headerLabels = @("165", "176", "219", "230", "255")

# New values for the 5 extra columns (E:I), rows 2-4, for every sheet.
$newData = @{
    "Sheet1" = @{
        2 = @(0.6767942583732057, 0.8659619528869733, 0.1046486298087396, 0, 0.0636084985984124)
        3 = @(0.2240525625420985, 0.8110732931652131, 0.07988905465319072, 0.3538156692436423, 0.6653323217601373)
        4 = @(0.1878588005374374, 0.8384741106513197, 0.4420361032894593, 0.4851463183440783, 0.7169509614167103)
    }
    "Sheet2" = @{
        2 = @(0.7045908183632734, 0.8538366619441524, 0.06291125020060985, 0, 0.0446162832744569)
        3 = @(0.207284754882557, 0.7782596172959202, 0.03953029831325686, 0.2034891273284664, 0.2670860790301472)
        4 = @(0.1227107772807832, 0.7151974233567935, 0.2360949233364733, 0.477882797731569, 0.544350779374509)
    }
    "Sheet3" = @{
        2 = @(0.2936982279542833, 0.06131674349908248, 0.0002579369344195344, 0, 0.01229528113913593)
        3 = @(0.00490418672236854, 0.2421089282605493, 0.005034612964128383, 0.02393376594316402, 0.1996670928521458)
        4 = @(0.01669144191972818, 0.006567216877041222, 0.001732220759328852, 0.1716231469871833, 0.2722694117579186)
    }
}

# Updated values for the pre-existing B:D columns, rows 2-4 (only Sheet2 and Sheet3 change).
$updatedBD = @{
    "Sheet2" = @{
        2 = @(0.5511376545814535, 0.7403166869671133, 0.8548771955004575)
        3 = @(0.3350805603868428, 0.2651338453581083, 0.2666733198369023)
        4 = @(0.09595175327723809, 0.4426125554850983, 0.06573705179282868)
    }
    "Sheet3" = @{
        2 = @(0.5554278761660133, 0.5503470475807263, 0.6047967180350287)
        3 = @(0.1029781286706216, 0.04475699175913309, 0.04596986541837332)
        4 = @(0.04456059601357413, 0.01417725635381229, 0.09742188428181062)
    }
}

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $name = $ws.Name

    # Update the pre-existing B2:D4 block where it changed.
    if ($updatedBD.ContainsKey($name)) {
        $rowsBD = $updatedBD[$name]
        foreach ($r in @(2, 3, 4)) {
            $vals = $rowsBD[$r]
            $ws.Range("B" + $r).Value = $vals[0]
            $ws.Range("C" + $r).Value = $vals[1]
            $ws.Range("D" + $r).Value = $vals[2]
        }
    }

    # Copy the header formatting (bold font, border, centered alignment) from
    # the existing B1:D1 header cells onto the new E1:I1 header cells before
    # writing their values, so the new headers share the same visual style.
    $ws.Range("B1:D1").Copy()
    $ws.Range("E1:I1").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # Write the five new header labels in row 1 (E1:I1), forcing text storage
    # (matching the existing B1:D1 text headers) via a leading quote-prefix.
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($headerCols[$i] + "1").Value = "'" + $headerLabels[$i]
    }

    # Write the new data values for E2:I4.
    $rowsNew = $newData[$name]
    foreach ($r in @(2, 3, 4)) {
        $vals = $rowsNew[$r]
        for ($i = 0; $i -lt 5; $i++) {
            $ws.Range($headerCols[$i] + $r).Value = $vals[$i]
        }
    }
}
